# Insert a new data row above row 79 (pushing existing rows 79-167 down to 80-168)
# and populate it with a new weekly price observation for "Repollo" at
# Terminal Hortofrutícola Agro Chillán.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 79; rows 79..167 shift down to 80..168.
$ws.Rows.Item(79).Insert()

# Seed the new row 79 with the same static field values as the row that is
# now directly below it (originally row 79, now row 80), since most columns
# (Mercado ID, Mercado, Región, Codreg, Categoría ID, Categoría, Variedad,
# Calidad, Unidad de comercialización, Origen, Kg o Unidades, Clasificación)
# are unchanged for this new observation.
$ws.Range("A80:R80").Copy()
$ws.Range("A79:R79").PasteSpecial()

# Now overwrite the fields that differ for this new observation
# (Fecha, Volumen, Precio mínimo, Precio máximo, Precio promedio ponderado,
# Precio $/Kg).
$ws.Range("D79").Value = 44554
$ws.Range("J79").Value = 1600
$ws.Range("K79").Value = 700
$ws.Range("L79").Value = 800
$ws.Range("M79").Value = 750
$ws.Range("P79").Value = 750
